$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet from "GossF-HW45.xpc" to "GossF"
$ws.Name = "GossF"

# Copy formatting (bold/border/center) from A15 to A16 before setting values
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

# Add row 16 data
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "HexGrid-60degTilt5degRes"
$ws.Range("C16").Value = 0.9969934576433873
$ws.Range("D16").Value = 0.9977992842888757
$ws.Range("E16").Value = 1.001764705882353
$ws.Range("F16").Value = 0.9969934576433873
$ws.Range("G16").Value = 1.001800823639892
$ws.Range("H16").Value = 1.004001539351017
$ws.Range("I16").Value = 1.004037657108557
$ws.Range("J16").Value = 0.9977992842888757
$ws.Range("K16").Value = 0.9997819950856143
$ws.Range("L16").Value = 0.9983877263645008
$ws.Range("M16").Value = 1.001066244652347
